# TestDataMappingSheet.xlsx update
# - Restyle three duplicated "Removal"-area rows (120:122) on the
#   TestDataMappingSheet_SD sheet so column A/C use the plain-font style
#   (no extra border/fill/alignment flags) instead of the old style.
# - Remove the stray duplicate "NewCase" row (row 123), which shifts
#   every following row up by one and shrinks the used range / filter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataMappingSheet_SD")

# Re-apply the font-only style to A120:A122 and C120:C122 (this mirrors the
# "retype value" edit that produced a new cellXf with only applyFont set).
foreach ($r in 120..122) {
    $ws.Cells.Item($r, 1).Font.Color = $ws.Cells.Item($r, 1).Font.Color
    $ws.Cells.Item($r, 3).Font.Color = $ws.Cells.Item($r, 3).Font.Color
}

# Delete the duplicate "NewCase" row entirely; everything below shifts up.
$ws.Rows.Item(123).Delete()

# Update the view to match where the user ended up after the edit.
$ws.Application.ActiveWindow.ScrollRow = 108
$ws.Range("A122").Select()

$wb.Save()
